$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row shifts left: A1 (班級) removed, B1/C1/D1 contents move into A1/B1/C1, D1 cleared
$ws.Range("A1").Value = "序號"
$ws.Range("B1").Value = "候選人"
$ws.Range("C1").Value = "得票數"
$ws.Range("D1").ClearContents()

# Update selection to match the committed state (D9)
$ws.Range("D9").Select()
